$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1089.8462
$ws.Range("I32").Value = 290.66666
$ws.Range("J32").Value = 1194.0869
$ws.Range("K32").Value = 290.66666
$ws.Range("L32").Value = 1194.0869
$ws.Range("M32").Value = 35.33334000000002
$ws.Range("N32").Value = -1846.0869
$ws.Range("H33").Value = 834.0526
$ws.Range("I33").Value = 856.7646999999999
$ws.Range("K33").Value = 856.7646999999999
$ws.Range("M33").Value = -627.7646999999999
$ws.Range("H121").Value = 1611.2
$ws.Range("J121").Value = 1683.4286
$ws.Range("L121").Value = 5050.2858
$ws.Range("N121").Value = -8544.2858
$ws.Range("H138").Value = 1633.8837
$ws.Range("I138").Value = 1218.4073
$ws.Range("J138").Value = 2335
$ws.Range("K138").Value = 3655.2219
$ws.Range("L138").Value = 7005
$ws.Range("M138").Value = 1484.7781
$ws.Range("N138").Value = -17285
$ws.Range("H141").Value = 6547.7856
$ws.Range("I141").Value = 2089.32
$ws.Range("K141").Value = 6267.960000000001
$ws.Range("M141").Value = -1087.960000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5221367.5
$ws.Range("I32").Value = 7574.093
$ws.Range("J32").Value = 33375852
$ws.Range("K32").Value = 7574.093
$ws.Range("L32").Value = 33375852
$ws.Range("M32").Value = -7287.093
$ws.Range("N32").Value = -33376426
$ws.Range("H45").Value = 2400.121
$ws.Range("I45").Value = 1266.55
$ws.Range("K45").Value = 1266.55
$ws.Range("M45").Value = -889.55
$ws.Range("H122").Value = 2062.577
$ws.Range("I122").Value = 1001
$ws.Range("J122").Value = 3301.0833
$ws.Range("K122").Value = 3003
$ws.Range("L122").Value = 9903.249899999999
$ws.Range("M122").Value = -553
$ws.Range("N122").Value = -14803.2499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2209.5
$ws.Range("I86").Value = 1737.8
$ws.Range("J86").Value = 2995.6667
$ws.Range("K86").Value = 1737.8
$ws.Range("L86").Value = 2995.6667
$ws.Range("M86").Value = -614.8
$ws.Range("N86").Value = -5241.6667
$ws.Range("H89").Value = 2209.5
$ws.Range("I89").Value = 1737.8
$ws.Range("J89").Value = 2995.6667
$ws.Range("K89").Value = 8689
$ws.Range("L89").Value = 14978.3335
$ws.Range("M89").Value = -3073
$ws.Range("N89").Value = -26210.3335
$ws.Range("H107").Value = 3802.95
$ws.Range("I107").Value = 4031.1667
$ws.Range("K107").Value = 4031.1667
$ws.Range("M107").Value = -2111.1667
$ws.Range("H134").Value = 5776.5405
$ws.Range("I134").Value = 2916.7856
$ws.Range("J134").Value = 7517.2607
$ws.Range("K134").Value = 8750.356800000001
$ws.Range("L134").Value = 22551.7821
$ws.Range("M134").Value = -6215.356800000001
$ws.Range("N134").Value = -27621.7821

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3079.7334
$ws.Range("I16").Value = 4520.1665
$ws.Range("J16").Value = 2119.4443
$ws.Range("K16").Value = 4520.1665
$ws.Range("L16").Value = 2119.4443
$ws.Range("M16").Value = -4233.1665
$ws.Range("N16").Value = -2693.4443
$ws.Range("H31").Value = 8931021
$ws.Range("I31").Value = 1812.7407
$ws.Range("J31").Value = 17244422
$ws.Range("K31").Value = 1812.7407
$ws.Range("L31").Value = 17244422
$ws.Range("M31").Value = -1517.7407
$ws.Range("N31").Value = -17245012
$ws.Range("H34").Value = 8931021
$ws.Range("I34").Value = 1812.7407
$ws.Range("J34").Value = 17244422
$ws.Range("K34").Value = 1812.7407
$ws.Range("L34").Value = 17244422
$ws.Range("M34").Value = -1610.7407
$ws.Range("N34").Value = -17244826
$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -9251
$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 28500
$ws.Range("L72").Value = 51000
$ws.Range("M72").Value = -26256
$ws.Range("H113").Value = 3079.7334
$ws.Range("I113").Value = 4520.1665
$ws.Range("J113").Value = 2119.4443
$ws.Range("K113").Value = 4520.1665
$ws.Range("L113").Value = 2119.4443
$ws.Range("M113").Value = -2350.1665
$ws.Range("N113").Value = -6459.4443
$ws.Range("H132").Value = 3909.3333
$ws.Range("I132").Value = 2793.5557
$ws.Range("J132").Value = 5583
$ws.Range("K132").Value = 8380.667099999999
$ws.Range("L132").Value = 16749
$ws.Range("M132").Value = -5850.667099999999
$ws.Range("N132").Value = -21809
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 469545.84
$ws.Range("I68").Value = 758.87036
$ws.Range("J68").Value = 1058255.1
$ws.Range("K68").Value = 2276.61108
$ws.Range("L68").Value = 3174765.3
$ws.Range("M68").Value = -1465.61108
$ws.Range("N68").Value = -3176387.3
$ws.Range("H71").Value = 469545.84
$ws.Range("I71").Value = 758.87036
$ws.Range("J71").Value = 1058255.1
$ws.Range("K71").Value = 6829.83324
$ws.Range("L71").Value = 9524295.9
$ws.Range("M71").Value = -2773.83324
$ws.Range("N71").Value = -9532407.9
$ws.Range("H98").Value = 9705.166999999999
$ws.Range("I98").Value = 560.8333
$ws.Range("J98").Value = 18849.5
$ws.Range("K98").Value = 1682.4999
$ws.Range("L98").Value = 56548.5
$ws.Range("M98").Value = -184.4999
$ws.Range("N98").Value = -59544.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1233275.2
$ws.Range("I102").Value = 1554086.6
$ws.Range("K102").Value = 1554086.6
$ws.Range("M102").Value = -1552464.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2002.2
$ws.Range("I22").Value = 2700.5
$ws.Range("J22").Value = 1827.625
$ws.Range("K22").Value = 2700.5
$ws.Range("L22").Value = 1827.625
$ws.Range("M22").Value = -2405.5
$ws.Range("N22").Value = -2417.625
$ws.Range("H27").Value = 2002.2
$ws.Range("I27").Value = 2700.5
$ws.Range("J27").Value = 1827.625
$ws.Range("K27").Value = 2700.5
$ws.Range("L27").Value = 1827.625
$ws.Range("M27").Value = -2593.5
$ws.Range("N27").Value = -2041.625
$ws.Range("H40").Value = 52633910
$ws.Range("I40").Value = 83334376
$ws.Range("J40").Value = 4543.4287
$ws.Range("K40").Value = 83334376
$ws.Range("L40").Value = 4543.4287
$ws.Range("M40").Value = -83334240
$ws.Range("N40").Value = -4815.4287
$ws.Range("H94").Value = 24165
$ws.Range("J94").Value = 24165
$ws.Range("L94").Value = 24165
$ws.Range("N94").Value = -25517
$ws.Range("H100").Value = 3322.6667
$ws.Range("I100").Value = 2905
$ws.Range("J100").Value = 3907.4
$ws.Range("K100").Value = 2905
$ws.Range("L100").Value = 3907.4
$ws.Range("M100").Value = -2364
$ws.Range("N100").Value = -4989.4
$ws.Range("H136").Value = 11906379
$ws.Range("I136").Value = 20001044
$ws.Range("J136").Value = 2462.0588
$ws.Range("K136").Value = 60003132
$ws.Range("L136").Value = 7386.176399999999
$ws.Range("M136").Value = -60000582
$ws.Range("N136").Value = -12486.1764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3104.56
$ws.Range("I132").Value = 2267.4666
$ws.Range("J132").Value = 4360.2
$ws.Range("K132").Value = 6802.399800000001
$ws.Range("L132").Value = 13080.6
$ws.Range("M132").Value = -4272.399800000001
$ws.Range("N132").Value = -18140.6
